$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$q = [char]34

# ---------------------------------------------------------------------------
# Row 22 - Hard Drive (full row, same shape as row 21)
# ---------------------------------------------------------------------------
$ws.Range("B21:P21").Copy()
$ws.Range("B22:P22").PasteSpecial(-4122)
$ws.Range("N22").Clear()

$ws.Range("B22").Formula = "=IF(E22=" + $q + $q + ",0,IF(F22=" + $q + $q + ",1,IF(H22=" + $q + $q + ",2,3)))"
$ws.Range("C22").Value = 17
$ws.Range("D22").Value = "Hard Drive"
$ws.Range("E22").Value = "250GB SATAII"
$ws.Range("F22").Value = "Hitachi"
$ws.Range("G22").Value = "HTS545025B9A300"
$ws.Range("H22").Value = "NewEgg"
$ws.Range("I22").Value = "N82E16822145255"
$ws.Range("J22").Value = "http://www.newegg.com/Product/Product.aspx?Item=N82E16822145255"
$ws.Hyperlinks.Add($ws.Range("J22"), "http://www.newegg.com/Product/Product.aspx?Item=N82E16822145255")
$ws.Range("J21").Copy()
$ws.Range("J22").PasteSpecial(-4122)
$ws.Range("K22").Value = 34.99
$ws.Range("L22").Value = 1
$ws.Range("M22").Formula = "=L22*`$L`$4"
$ws.Range("O22").Formula = "=M22*K22+N22"
$ws.Range("P22").Value = "Discount ends 8/25"

# ---------------------------------------------------------------------------
# Rows 23-26 - misc. on-board connectors (sparse rows: B,C,D,K,L,M,O only)
# ---------------------------------------------------------------------------
$rows = @(
  @{ Row = 23; Item = 18; Name = "CPU Fan" },
  @{ Row = 24; Item = 19; Name = "24 pin Power Connector" },
  @{ Row = 25; Item = 20; Name = "12v Power Connector" },
  @{ Row = 26; Item = 21; Name = "Fan Connectors" }
)

foreach ($r in $rows) {
  $n = $r.Row

  $ws.Range("B21").Copy(); $ws.Range("B$n").PasteSpecial(-4122)
  $ws.Range("C21").Copy(); $ws.Range("C$n").PasteSpecial(-4122)
  $ws.Range("D21").Copy(); $ws.Range("D$n").PasteSpecial(-4122)
  $ws.Range("K21").Copy(); $ws.Range("K$n").PasteSpecial(-4122)
  $ws.Range("L21").Copy(); $ws.Range("L$n").PasteSpecial(-4122)
  $ws.Range("M21").Copy(); $ws.Range("M$n").PasteSpecial(-4122)
  $ws.Range("O21").Copy(); $ws.Range("O$n").PasteSpecial(-4122)

  $ws.Range("B$n").Formula = "=IF(E$n=" + $q + $q + ",0,IF(F$n=" + $q + $q + ",1,IF(H$n=" + $q + $q + ",2,3)))"
  $ws.Range("C$n").Value = $r.Item
  $ws.Range("D$n").Value = $r.Name
  $ws.Range("L$n").Value = 1
  $ws.Range("M$n").Formula = "=L$n*`$L`$4"
  $ws.Range("O$n").Formula = "=M$n*K$n+N$n"
}

# ---------------------------------------------------------------------------
# Named range pointing at the discounted hard-drive row
# ---------------------------------------------------------------------------
$ws.Names.Add("CART_ITEM", "=Sheet1!`$G`$22")

# ---------------------------------------------------------------------------
# View state - scrolled / selection moved while entering the new rows
# ---------------------------------------------------------------------------
$ws.Range("N26").Select()
try {
  $excel.ActiveWindow.ScrollColumn = 7
  $excel.ActiveWindow.ScrollRow = 4
} catch {
}
